$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price/Volume columns to remain plain text so values such as
# "1.006" or "29.537.44" are not reinterpreted as numbers/dates by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

$rows = @(
    @{Row=2; B="Bitcoin"; C="https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc"; D="29.537.44"; E="  +1.05%  "},
    @{Row=3; B="Ethereum"; C="https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth"; D="1.980.35"; E="  +4.20%  "},
    @{Row=4; B="TetherUSD"; C="https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt"; D="1.006"; E="  +0.43%  "},
    @{Row=5; B="BNB"; C="https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"; D="327.73"; E="  +0.36%  "},
    @{Row=6; B="USDC"; C="https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"; D="1.005"; E="  +0.39%  "},
    @{Row=7; B="XRP"; C="https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"; D="0.4663"; E="  +0.20%  "},
    @{Row=8; B="Cardano"; C="https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"; D="0.3916"; E="  -0.09%  "},
    @{Row=9; B="OKB"; C="https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"; D="46.25"; E="  -1.33%  "},
    @{Row=10; B="Dogecoin"; C="https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"; D="0.07958"; E="  +0.91%  "},
    @{Row=11; B="Polygon"; C="https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"; D="0.9932"; E="  +0.50%  "},
    @{Row=12; B="Solana"; C="https://coinranking.com/coin/zNZHO_Sjf+solana-sol"; D="22.93"; E="  +4.30%  "},
    @{Row=13; B="WrappedEther"; C="https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"; D="2.004.84"; E="  +4.69%  "},
    @{Row=14; B="Chainlink"; C="https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"; D="7.191"; E="  +1.60%  "},
    @{Row=15; B="Polkadot"; C="https://coinranking.com/coin/25W7FG7om+polkadot-dot"; D="5.838"; E="  +1.52%  "},
    @{Row=16; B="TRON"; C="https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"; D="0.07105"; E="  +1.93%  "},
    @{Row=17; B="Litecoin"; C="https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"; D="87.78"; E="  -0.66%  "},
    @{Row=18; B="BinanceUSD"; C="https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"; D="1.007"; E="  +0.41%  "},
    @{Row=19; B="ShibaInu"; C="https://coinranking.com/coin/xz24e0BjL+shibainu-shib"; D="0.000009974"; E="  -0.08%  "},
    @{Row=20; B="Avalanche"; C="https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"; D="17.30"; E="  +1.26%  "},
    @{Row=21; B="Dai"; C="https://coinranking.com/coin/MoTuySvg7+dai-dai"; D="1.004"; E="  +0.35%  "},
    @{Row=22; B="WrappedBTC"; C="https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"; D="29.555.20"; E="  +1.04%  "},
    @{Row=23; B="BitDAO"; C="https://coinranking.com/coin/N2IgQ9Xme+bitdao-bit"; D="0.5112"; E="  +6.75%  "},
    @{Row=24; B="Uniswap"; C="https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"; D="5.562"; E="  +4.61%  "},
    @{Row=25; B="Cosmos"; C="https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"; D="11.20"; E="  +1.00%  "},
    @{Row=26; B="WrappedliquidstakedEther2.0"; C="https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"; D="2.235.88"; E="  +4.18%  "},
    @{Row=27; B="Toncoin"; C="https://coinranking.com/coin/67YlI0K1b+toncoin-ton"; D="2.111"; E="  +0.75%  "},
    @{Row=28; B="Monero"; C="https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"; D="158.51"; E="  +1.42%  "},
    @{Row=29; B="EthereumClassic"; C="https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"; D="19.64"; E="  +0.94%  "},
    @{Row=30; B="InternetComputer(DFINITY)"; C="https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"; D="5.831"; E="  -2.54%  "},
    @{Row=31; B="BitcoinCash"; C="https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"; D="119.65"; E="  +0.85%  "},
    @{Row=32; B="LidoDAOToken"; C="https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"; D="1.902"; E="  -0.39%  "},
    @{Row=33; B="Stellar"; C="https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"; D="0.09434"; E="  +0.82%  "},
    @{Row=34; B="ImmutableX"; C="https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"; D="0.8922"; E="  -1.55%  "},
    @{Row=35; B="Filecoin"; C="https://coinranking.com/coin/ymQub4fuB+filecoin-fil"; D="5.238"; E="  -0.90%  "},
    @{Row=36; B="ARBITRUM"; C="https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"; D="1.325"; E="  -0.01%  "},
    @{Row=37; B="HuobiToken"; C="https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"; D="3.200"; E="  -0.45%  "},
    @{Row=38; B="Hedera"; C="https://coinranking.com/coin/jad286TjB+hedera-hbar"; D="0.05818"; E="  +0.61%  "},
    @{Row=39; B="TrustWalletToken"; C="https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"; D="1.179"; E="  -0.36%  "},
    @{Row=40; B="VeChain"; C="https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"; D="0.02100"; E="  +0.53%  "},
    @{Row=41; B="FraxShare"; C="https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"; D="7.803"; E="  +0.72%  "},
    @{Row=42; B="TheSandbox"; C="https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"; D="0.5738"; E="  +0.47%  "},
    @{Row=43; B="Algorand"; C="https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"; D="0.1806"; E="  +1.05%  "},
    @{Row=44; B="PEPE"; C="https://coinranking.com/coin/03WI8NQPF+pepe-pepe"; D="0.000003037"; E="  +35.24%  "},
    @{Row=45; B="Aptos"; C="https://coinranking.com/coin/HGYj5JCv5+aptos-apt"; D="9.695"; E="  -0.61%  "},
    @{Row=46; B="MXToken"; C="https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"; D="2.788"; E="  +8.18%  "},
    @{Row=47; B="Decentraland"; C="https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"; D="0.5370"; E="  +0.37%  "},
    @{Row=48; B="EnergySwap"; C="https://coinranking.com/coin/SbWqqTui-+energyswap-ens"; D="11.77"; E="  -1.50%  "},
    @{Row=49; B="RenderToken"; C="https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"; D="2.160"; E="  -1.64%  "},
    @{Row=50; B="Cronos"; C="https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"; D="0.06939"; E="  -1.52%  "},
    @{Row=51; B="Quant"; C="https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"; D="114.24"; E="  +0.86%  "},
)

foreach ($item in $rows) {
    $ws.Cells.Item($item.Row, 2).Value = $item.B
    $ws.Cells.Item($item.Row, 3).Value = $item.C
    $ws.Cells.Item($item.Row, 4).Value = $item.D
    $ws.Cells.Item($item.Row, 5).Value = $item.E
}
